# "further progress with Newsletter"
#
# Adds a 5th author column ("author5") to the Newsletter sheet: a new
# header in T1 and a default value of "Anonym" for every existing data
# row (T2:T6), mirroring the pattern already used by author1..author4
# (columns O:R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column T
$ws.Range("T1").Value = "author5"

# Default author value for the five existing article rows
$ws.Range("T2:T6").Value = "Anonym"

# Leave the selection where the author ended up after entering the data
$ws.Range("T11").Select()
